# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the refreshed counts recorded at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    switch ($sheetName) {
        "展览" {
            $updates = @{
                3  = 546
                4  = 1566
                12 = 355
                14 = 510
                16 = 6521
                17 = 28
                22 = 15579
                23 = 1542
                27 = 11119
                28 = 779
                29 = 4361
                30 = 255
                34 = 130
            }
        }
        "全部类型" {
            $updates = @{
                3  = 546
                4  = 1566
                14 = 355
                16 = 510
                19 = 6521
                20 = 28
                26 = 15579
                27 = 1542
                32 = 11119
                33 = 779
                34 = 4361
                35 = 255
                39 = 130
            }
        }
    }

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
